# Updated cryptos list on Mon Dec 25 17:35:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the Price column (D) stores values as literal text in the source
# workbook (e.g. "7.00", "1.00", thousand-dot formats like "43.773.94").
# Many of the new prices are plain decimals that Excel's COM layer would
# otherwise auto-coerce to a Number (dropping trailing zeros / introducing
# float noise), so for those cells we force the Text number format first to
# preserve the exact string, matching the source data's inline-string type.

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "43.777.47"
$ws.Range("E2").Value = "  -0.02%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.281.76"
$ws.Range("E3").Value = "  -0.40%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.04%  "

# --- Row 5: Solana ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "123.73"
$ws.Range("E5").Value = "  +9.04%  "

# --- Row 6: BNB ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.52"
$ws.Range("E6").Value = "  -0.68%  "

# --- Row 7: XRP ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +2.02%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.23%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +1.36%  "

# --- Row 10: Avalanche ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.47"
$ws.Range("E10").Value = "  -0.11%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  +0.69%  "

# --- Row 12: Polkadot ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.21"
$ws.Range("E12").Value = "  +2.00%  "

# --- Row 13: TRON ---
$ws.Range("E13").Value = "  -0.99%  "

# --- Row 14: Chainlink ---
$ws.Range("E14").Value = "  -1.85%  "

# --- Row 15: Polygon ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.907"
$ws.Range("E15").Value = "  +5.07%  "

# --- Row 16: WrappedliquidstakedEther2.0 ---
$ws.Range("D16").Value = "2.624.11"
$ws.Range("E16").Value = "  -0.49%  "

# --- Row 17: WrappedEther ---
$ws.Range("D17").Value = "2.280.04"
$ws.Range("E17").Value = "  -0.84%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "43.716.01"
$ws.Range("E18").Value = "  +0.01%  "

# --- Row 19: ShibaInu ---
$ws.Range("E19").Value = "  +0.66%  "

# --- Row 20: Uniswap ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  -0.09%  "

# --- Row 21: Litecoin ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.44"
$ws.Range("E21").Value = "  +0.22%  "

# --- Row 22: ImmutableX ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").Value = "  +1.09%  "

# --- Row 23: BitcoinCash ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.35"
$ws.Range("E23").Value = "  +1.48%  "

# --- Row 24 & 25: PancakeSwap / InternetComputer(DFINITY) swap places ---
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.46"
$ws.Range("E25").Value = "  -4.09%  "

# --- Row 26: Dai ---
$ws.Range("E26").Value = "  +1.91%  "

# --- Row 27: Cosmos ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.85"
$ws.Range("E27").Value = "  +1.03%  "

# --- Row 28: InjectiveProtocol ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.94"
$ws.Range("E28").Value = "  -0.74%  "

# --- Row 29: WEMIXToken ---
$ws.Range("E29").Value = "  -0.53%  "

# --- Row 30: Toncoin ---
$ws.Range("E30").Value = "  +0.25%  "

# --- Row 31: Monero ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.65"
$ws.Range("E31").Value = "  -1.06%  "

# --- Row 32: EthereumClassic ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.75"
$ws.Range("E32").Value = "  +0.69%  "

# --- Row 33: Hedera ---
$ws.Range("E33").Value = "  -1.75%  "

# --- Row 34: Filecoin ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.76"
$ws.Range("E34").Value = "  +1.49%  "

# --- Row 35: Stellar ---
$ws.Range("E35").Value = "  +1.99%  "

# --- Row 36: NEARProtocol ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.28"
$ws.Range("E36").Value = "  +12.41%  "

# --- Row 37: VeChain ---
$ws.Range("E37").Value = "  +5.22%  "

# --- Row 38: RenderToken ---
$ws.Range("E38").Value = "  -2.39%  "

# --- Row 39: Kaspa ---
$ws.Range("E39").Value = "  +0.85%  "

# --- Row 40: LidoDAOToken ---
$ws.Range("E40").Value = "  +5.85%  "

# --- Row 41: MultiversX ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.27"
$ws.Range("E41").Value = "  +1.00%  "

# --- Row 42: Celestia ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.02"
$ws.Range("E42").Value = "  -4.13%  "

# --- Row 43: Algorand ---
$ws.Range("E43").Value = "  -0.62%  "

# --- Row 44: FirstDigitalUSD ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.19%  "

# --- Row 45: ARBITRUM ---
$ws.Range("E45").Value = "  -2.57%  "

# --- Row 46: THORChain ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.66"
$ws.Range("E46").Value = "  -10.56%  "

# --- Row 47: ordi ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.88"
$ws.Range("E47").Value = "  +37.00%  "

# --- Row 48: FraxShare ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.63"
$ws.Range("E48").Value = "  -1.80%  "

# --- Row 49: TrustWalletToken ---
$ws.Range("E49").Value = "  +0.37%  "

# --- Row 50: Cronos ---
$ws.Range("E50").Value = "  +0.91%  "

# --- Row 51: Aave ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.24"
$ws.Range("E51").Value = "  -0.61%  "
